# Add a second "published_Colloff_SealeCarlisle_Karoğlu_etal2020_E2()" block
# of model-fit results (12 new columns, AL:AW, covering 2 conditions x 6
# parameter-set columns each) to the right of the existing E1() block
# (columns B:AK, 6 conditions x 6 columns each), mirroring rows 1-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCols = @("AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW")

# xlPasteFormats - used to stamp the correct cell style after typing
# brand-new text values (typed values otherwise pick up the plain
# column-default style instead of the header/bold style used nearby).
$xlPasteFormats = -4122

# --- Row 1: new experiment header ------------------------------------
foreach ($c in $targetCols) {
    $ws.Range($c + "1").Value = "published_Colloff_SealeCarlisle_Karoğlu_etal2020_E2()"
}
$ws.Range("AK1").Copy()
$ws.Range("AL1:AW1").PasteSpecial($xlPasteFormats)

# --- Row 2: Model fit (IndependentObservation / Ensemble / Integration)
$row2Src = @("B","B","D","D","F","F","B","B","D","D","F","F")
for ($i = 0; $i -lt $targetCols.Length; $i++) {
    $ws.Range($row2Src[$i] + "2").Copy($ws.Range($targetCols[$i] + "2"))
}

# --- Row 3: Parameter set (EqualVariance / UnequalVariance) ----------
$row3Src = @("B","C","B","C","B","C","B","C","B","C","B","C")
for ($i = 0; $i -lt $targetCols.Length; $i++) {
    $ws.Range($row3Src[$i] + "3").Copy($ws.Range($targetCols[$i] + "3"))
}

# --- Row 4: new exclusions spec ({}) ----------------------------------
foreach ($c in $targetCols) {
    $ws.Range($c + "4").Value = "{}"
}
$ws.Range("AK4").Copy()
$ws.Range("AL4:AW4").PasteSpecial($xlPasteFormats)

# --- Row 5: Condition (condition 1 / condition 2) ---------------------
$row5Src = @("B","B","B","B","B","B","H","H","H","H","H","H")
for ($i = 0; $i -lt $targetCols.Length; $i++) {
    $ws.Range($row5Src[$i] + "5").Copy($ws.Range($targetCols[$i] + "5"))
}

# --- Row 6: Binning ([-1,60,80,100]), uniform --------------------------
foreach ($c in $targetCols) {
    $ws.Range("AK6").Copy($ws.Range($c + "6"))
}

# --- Row 7: Exclusions (True), uniform ---------------------------------
foreach ($c in $targetCols) {
    $ws.Range("AK7").Copy($ws.Range($c + "7"))
}

# --- Row 8: niter = 2000, uniform ---------------------------------------
foreach ($c in $targetCols) {
    $ws.Range("AK8").Copy($ws.Range($c + "8"))
}

# --- Column widths for the new block (mirrors the B:AK pattern) --------
$ws.Range("AL1:AL43").ColumnWidth = 18.666666666666668
$ws.Range("AM1:AQ43").ColumnWidth = 16
$ws.Range("AR1:AR43").ColumnWidth = 18.666666666666668
$ws.Range("AS1:AW43").ColumnWidth = 16

# --- Selection moves to the newly-added block --------------------------
$ws.Range("AL11").Select()

"done"
